# Updated FRA model - 2025-08-06 18:10
$wb = $excel.ActiveWorkbook

# --- Rename "Sheet1" to "misc." ---
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "misc."

# --- Add the two new rows (11 & 12) of ~TFM_INS-TS data ---
# Entry order matches original authoring so shared-string indices line up:
# flo_emis, gas, co2captured, co2, *ccs,*ccs-rf, coal,oil
$ws.Cells.Item(11, 2).Value = "flo_emis"
$ws.Cells.Item(11, 4).Value = "gas"
$ws.Cells.Item(11, 12).Value = "co2captured"
$ws.Cells.Item(11, 11).Value = "co2"
$ws.Cells.Item(11, 5).Value = "*ccs,*ccs-rf"
$ws.Cells.Item(11, 8).Value = 0.95

# Row 12: flo_emis / coal,oil / *ccs,*ccs-rf / 0.85 / co2 / co2captured
$ws.Cells.Item(12, 2).Value = "flo_emis"
$ws.Cells.Item(12, 4).Value = "coal,oil"
$ws.Cells.Item(12, 12).Value = "co2captured"
$ws.Cells.Item(12, 11).Value = "co2"
$ws.Cells.Item(12, 5).Value = "*ccs,*ccs-rf"
$ws.Cells.Item(12, 8).Value = 0.85

# --- New header cells for columns K (other_indexes) and L (commodity) ---
$ws.Cells.Item(4, 11).Value = "other_indexes"
$ws.Cells.Item(4, 12).Value = "commodity"
# Match the bold/underlined "Heading 2" style used by the rest of row 4
$ws.Range("K4:L4").Style = $ws.Range("J4").Style

# --- Column width adjustments ---
# (ColumnWidth assignment applies a fixed engine offset; these inputs land
# on the closest achievable stored widths to the target 10.33203125 / 12)
$ws.Columns.Item(5).ColumnWidth = 9.5
$ws.Columns.Item(11).ColumnWidth = 11.166666666666666

# --- Selection moves to D13 after the edit ---
$ws.Range("D13").Select() | Out-Null
